# Dragons Den Poster - "Final Changes to Presentation Notes"
#
# Maps the author's edits onto the single-slide poster deck. Shapes are
# addressed by their stable z-order index within Slide 1:
#   4  -> id 6  "Rectangle 5"   (was blank)
#   5  -> id 7  "Rectangle 6"   ("SIMPLICITY IN NAVIGATION")
#   6  -> id 8  "Rectangle 7"   ("WHAT CAN MEMORI DO" + extra paragraphs)
#   7  -> id 9  "Rectangle 8"   (was blank)
#   8  -> id 11 "Rectangle 10"  ("WHY CHOOSE MEMORI")
#   9  -> id 12 "Rectangle 11"  ("LIGHT/DARK MODE")
#   10 -> id 14 "Rectangle 13"  ("WHAT IS MEMORI" + extra paragraphs)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------
# "Rectangle 5" (blank) gains the heading "WHAT CAN MEMORI DO?"
# ---------------------------------------------------------------
$sh = $s.Shapes.Item(4)
$tr = $sh.TextFrame.TextRange
$tr.Text = "WHAT CAN MEMORI DO?"
$tr.LanguageID = "en-GB"

# ---------------------------------------------------------------
# "Rectangle 6" moves/resizes and its caption becomes "IMAGE?"
# ---------------------------------------------------------------
$sh = $s.Shapes.Item(5)
$sh.Left = 0
$sh.Top = 371.5
$sh.Width = 209.5
$sh.Height = 168.5
$tr = $sh.TextFrame.TextRange
$tr.Text = "IMAGE?"
$tr.LanguageID = "en-GB"

# ---------------------------------------------------------------
# "Rectangle 7" collapses its three paragraphs into one "IMAGE?"
# ---------------------------------------------------------------
$sh = $s.Shapes.Item(6)
$tr = $sh.TextFrame.TextRange
$tr.Text = "IMAGE?"
$tr.LanguageID = "en-GB"

# ---------------------------------------------------------------
# "Rectangle 8" (blank) gains the heading "WHAT MAKES MEMORI THE BEST?"
# ---------------------------------------------------------------
$sh = $s.Shapes.Item(7)
$tr = $sh.TextFrame.TextRange
$tr.Text = "WHAT MAKES MEMORI THE BEST?"
$tr.LanguageID = "en-GB"

# ---------------------------------------------------------------
# "Rectangle 10" caption becomes "IMAGE?"
# ---------------------------------------------------------------
$sh = $s.Shapes.Item(8)
$tr = $sh.TextFrame.TextRange
$tr.Text = "IMAGE?"
$tr.LanguageID = "en-GB"

# ---------------------------------------------------------------
# "Rectangle 11" moves/resizes and its caption becomes "WHY CHOOSE MEMORI?"
# ---------------------------------------------------------------
$sh = $s.Shapes.Item(9)
$sh.Left = 230.5
$sh.Top = 361
$sh.Width = 224.0001
$sh.Height = 179
$tr = $sh.TextFrame.TextRange
$tr.Text = "WHY CHOOSE MEMORI?"
$tr.LanguageID = "en-GB"

# ---------------------------------------------------------------
# "Rectangle 13" collapses its three paragraphs into one, split across
# two runs: "WHAT IS " + "MEMORI"
# ---------------------------------------------------------------
$sh = $s.Shapes.Item(10)
$tr = $sh.TextFrame.TextRange
$tr.Text = "WHAT IS "
$tr.LanguageID = "en-GB"
$tr2 = $tr.InsertAfter("MEMORI")
$tr2.LanguageID = "en-GB"
